$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connectors")

# The digital connector labels in row 3 were off by one (reading sensor /
# labeling connectors bug). Renumber D1..D6 -> D2..D7.
$ws.Range("L3").Value = "D7"
$ws.Range("K3").Value = "D6"
$ws.Range("J3").Value = "D5"
$ws.Range("I3").Value = "D4"
$ws.Range("H3").Value = "D3"
$ws.Range("G3").Value = "D2"

# Reflect the active selection on the sheet as it was when the fix was made
$ws.Range("L3").Select()
